# DatabaseSamples.xlsx update
# Adds a new student (Dristi Marasini, s533985 / id 533985) to the
# "Student" sheet and the corresponding term history rows to the
# "StudentTerm" sheet.

$wb = $excel.ActiveWorkbook

# --- Student sheet: add new student record (row 3) ---
$wsStudent = $wb.Worksheets.Item("Student")
$wsStudent.Range("A3").Value = 533985
$wsStudent.Range("B3").Value = "Dristi"
$wsStudent.Range("C3").Value = "Marasini"
$wsStudent.Range("D3").Value = "s533985"
$wsStudent.Range("E3").Value = 569178
$wsStudent.Range("E3").Select()

# --- StudentTerm sheet: add term history rows 8-12 for the new student ---
$wsStudentTerm = $wb.Worksheets.Item("StudentTerm")

$wsStudentTerm.Range("A8").Value = 7
$wsStudentTerm.Range("B8").Value = 533985
$wsStudentTerm.Range("C8").Value = 7
$wsStudentTerm.Range("D8").Value = "Fall 2018"

$wsStudentTerm.Range("A9").Value = 8
$wsStudentTerm.Range("B9").Value = 533985
$wsStudentTerm.Range("C9").Value = 8
$wsStudentTerm.Range("D9").Value = "Spring 2019"

$wsStudentTerm.Range("A10").Value = 9
$wsStudentTerm.Range("B10").Value = 533985
$wsStudentTerm.Range("C10").Value = 9
$wsStudentTerm.Range("D10").Value = "Summer 2019"

$wsStudentTerm.Range("A11").Value = 10
$wsStudentTerm.Range("B11").Value = 533985
$wsStudentTerm.Range("C11").Value = 10
$wsStudentTerm.Range("D11").Value = "Fall 2019"

$wsStudentTerm.Range("A12").Value = 11
$wsStudentTerm.Range("B12").Value = 533985
$wsStudentTerm.Range("C12").Value = 11
$wsStudentTerm.Range("D12").Value = "Spring 2019"

$wsStudentTerm.Range("C13").Select()
